$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.121.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.508.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '109.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '320.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("E10").Value = '  +3.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.32%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.124'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.902.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.512.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.960.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("E22").Value = '  +2.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '273.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.56%  '
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("E28").Value = '  +5.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.142'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("E33").Value = '  -5.79%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0783'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '121.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.26%  '
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0309'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.031.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.86%  '
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.08%  '
